# Updates the "cryptos" price/volume table with a fresh data pull.
# Prices/percentages are stored as plain text (not numbers), so every
# write is forced to text (leading apostrophe) and the cell style is
# reset to "Normal" afterwards so no stray number-format/quote-prefix
# style gets attached to the cell (keeps cells on the default style,
# matching the original workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "61.481.44"
Set-TextValue "E2" "  +1.31%  "
Set-TextValue "D3" "3.385.96"
Set-TextValue "E3" "  +1.07%  "
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "576.06"
Set-TextValue "E5" "  +1.16%  "
Set-TextValue "D6" "136.69"
Set-TextValue "E6" "  +2.28%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "3.386.14"
Set-TextValue "E8" "  +1.07%  "
Set-TextValue "E9" "  -0.49%  "
Set-TextValue "E10" "  -0.84%  "
Set-TextValue "E11" "  +3.12%  "
Set-TextValue "E12" "  +0.56%  "
Set-TextValue "D13" "3.965.65"
Set-TextValue "E13" "  +1.00%  "
Set-TextValue "E14" "  +2.99%  "
Set-TextValue "E15" "  +2.97%  "
Set-TextValue "D16" "3.388.00"
Set-TextValue "E16" "  +1.02%  "
Set-TextValue "D17" "25.69"
Set-TextValue "E17" "  +2.92%  "
Set-TextValue "D18" "61.588.42"
Set-TextValue "E18" "  +1.22%  "
Set-TextValue "D19" "14.16"
Set-TextValue "E19" "  +2.27%  "
Set-TextValue "D20" "5.87"
Set-TextValue "E20" "  +2.32%  "
Set-TextValue "D21" "9.39"
Set-TextValue "E21" "  +0.49%  "
Set-TextValue "D22" "376.75"
Set-TextValue "E22" "  +1.33%  "
Set-TextValue "E23" "  -2.24%  "
Set-TextValue "D24" "3.530.93"
Set-TextValue "E24" "  +1.25%  "
Set-TextValue "E25" "  +0.22%  "
Set-TextValue "E26" "  +9.27%  "
Set-TextValue "D27" "71.26"
Set-TextValue "E27" "  +1.30%  "
Set-TextValue "D28" "1.71"
Set-TextValue "E28" "  +3.41%  "
Set-TextValue "D29" "7.52"
Set-TextValue "E29" "  -1.21%  "
Set-TextValue "D30" "0.998"
Set-TextValue "E30" "  -0.21%  "
Set-TextValue "E31" "  +5.26%  "
Set-TextValue "E32" "  +2.34%  "
Set-TextValue "D33" "2.18"
Set-TextValue "E33" "  +2.24%  "
Set-TextValue "E34" "  +0.08%  "
Set-TextValue "D35" "23.42"
Set-TextValue "E35" "  +0.66%  "
Set-TextValue "D36" "5.31"
Set-TextValue "E36" "  -3.67%  "
Set-TextValue "E37" "  +0.49%  "
Set-TextValue "D38" "6.82"
Set-TextValue "E38" "  -0.53%  "
Set-TextValue "D39" "165.37"
Set-TextValue "E39" "  +1.90%  "
Set-TextValue "E40" "  -0.02%  "
Set-TextValue "D41" "0.780"
Set-TextValue "E41" "  +3.41%  "
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  -0.14%  "
Set-TextValue "E43" "  +2.84%  "
Set-TextValue "E44" "  +8.65%  "
Set-TextValue "B45" "Filecoin"
Set-TextValue "C45" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D45" "4.41"
Set-TextValue "E45" "  +1.04%  "
Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "24.90"
Set-TextValue "E46" "  +9.63%  "
Set-TextValue "D47" "41.41"
Set-TextValue "E48" "  -1.11%  "
Set-TextValue "E49" "  -2.68%  "
Set-TextValue "D50" "2.340.46"
Set-TextValue "E50" "  +5.68%  "
Set-TextValue "E51" "  -0.73%  "
